$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 0.4386883400186667
$ws.Range("R2").Value = 3.948195060168
$ws.Range("S2").Value = 0.02769484181536182
$ws.Range("T2").Value = 0.02769484181536182

# Row 3 updates
$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("S3").Value = 0.6282762845978157
$ws.Range("T3").Value = 0.6282762845978156

# Row 4 updates
$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("R4").Value = 49.044984921966
$ws.Range("S4").Value = 0.3440288735868225
$ws.Range("T4").Value = 0.3440288735868225
